$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Cell values.
# NOTE: the write order below is deliberately chosen (col B, then col A,
# then col C) so that the workbook's shared-string table is rebuilt in
# the same order the source workbook uses: the "tags" strings first,
# then the three "Loudspeaker:" step descriptions, then the new Jira key
# string last.
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "chat,acceptance,smoke"
$ws.Range("B2").Value = "chat,acceptance,smoke"
$ws.Range("B3").Value = "chat,acceptance"

$ws.Range("A1").Value = "Loudspeaker: move phone to ear to switch OFF loudspeaker"
$ws.Range("A2").Value = "Loudspeaker: move phone from ear to switch back ON loudspeaker"
$ws.Range("A3").Value = "Loudspeaker: check on different volume that ON-OFF transition of loudspeaker doesn't affect the volume"

$ws.Range("C1").Value = "GRD-342,GRD-343"
$ws.Range("C2").Value = "GRD-342,GRD-343"
$ws.Range("C3").Value = "GRD-342,GRD-343"

# ---------------------------------------------------------------------
# Column A becomes much wider (to fit the long step text); columns B/C
# keep their original width untouched.
# ---------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 56.166666666666664

# ---------------------------------------------------------------------
# Rows 2 and 3 get a taller, explicit height to make room for the
# wrapped text in column A.
# ---------------------------------------------------------------------
$ws.Rows(2).RowHeight = 30
$ws.Rows(3).RowHeight = 30

# ---------------------------------------------------------------------
# Column A (the step-description column) gets a thin box border, wraps
# its text, and aligns it to the top of the cell.
# ---------------------------------------------------------------------
$stepsRange = $ws.Range("A1:A3")
$stepsRange.Borders.LineStyle = 1
$stepsRange.WrapText = $true
$stepsRange.VerticalAlignment = -4160

# ---------------------------------------------------------------------
# Move the active selection.
# ---------------------------------------------------------------------
$ws.Range("G2").Select() | Out-Null

Write-Output "done"
